$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.813.25"
$ws.Range("E2").Value = "  -2.50%  "

$ws.Range("D3").Value = "2.913.30"
$ws.Range("E3").Value = "  -3.34%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "584.74"
$ws.Range("E5").Value = "  -1.59%  "

$ws.Range("D6").Value = "146.65"
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("D8").Value = "2.908.14"
$ws.Range("E8").Value = "  -3.46%  "

$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +7.45%  "

$ws.Range("D11").Value = "0.143"
$ws.Range("E11").Value = "  -3.78%  "

$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -3.64%  "

$ws.Range("D14").Value = "34.31"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").Value = "3.395.85"
$ws.Range("E16").Value = "  -3.34%  "

$ws.Range("D17").Value = "60.836.41"
$ws.Range("E17").Value = "  -2.40%  "

$ws.Range("D18").Value = "6.80"
$ws.Range("E18").Value = "  -2.62%  "

$ws.Range("D19").Value = "2.916.00"
$ws.Range("E19").Value = "  -3.31%  "

$ws.Range("D20").Value = "424.89"
$ws.Range("E20").Value = "  -5.48%  "

$ws.Range("D21").Value = "13.65"
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("D22").Value = "0.669"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").Value = "7.16"
$ws.Range("E23").Value = "  -2.79%  "

$ws.Range("D24").Value = "80.75"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").Value = "10.98"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -3.55%  "

$ws.Range("D27").Value = "11.81"
$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").Value = "2.62"
$ws.Range("E31").Value = "  -3.13%  "

$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  +3.23%  "

$ws.Range("D33").Value = "26.66"
$ws.Range("E33").Value = "  -2.64%  "

$ws.Range("E34").Value = "  -4.04%  "

$ws.Range("D35").Value = "0.0₃0841"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("E36").Value = "  -1.67%  "

$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  -2.90%  "

$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("D39").Value = "49.78"
$ws.Range("E39").Value = "  -1.04%  "

$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("D41").Value = "0.123"
$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("D42").Value = "8.76"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("D43").Value = "41.93"
$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("D44").Value = "0.287"
$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("D45").Value = "373.76"
$ws.Range("E45").Value = "  -7.52%  "

$ws.Range("D46").Value = "0.0345"
$ws.Range("E46").Value = "  -2.04%  "

$ws.Range("D47").Value = "2.650.78"
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("D48").Value = "132.85"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").Value = "25.42"
$ws.Range("E49").Value = "  +7.02%  "

$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("E51").Value = "  -1.10%  "
